# Weekly update: insert a new data row for "Feria Lagunitas de Puerto Montt" /
# "Cebollín" right before the current row 210, shifting the existing rows
# 210..224 down to 211..225, and populate the newly inserted row with the
# latest week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 210 (pushes old 210..224 down to 211..225)
$ws.Rows.Item(210).Insert()

# Populate the newly inserted row 210 with the new weekly record
$ws.Range("A210").Value2 = 4
$ws.Range("B210").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C210").Value2 = "Los Lagos"
$ws.Range("D210").Value2 = 44585
$ws.Range("E210").Value2 = 10
$ws.Range("F210").Value2 = 100112037
$ws.Range("G210").Value2 = "Cebollín"
$ws.Range("H210").Value2 = "Sin especificar"
$ws.Range("I210").Value2 = "Primera"
$ws.Range("J210").Value2 = 70
$ws.Range("K210").Value2 = 6000
$ws.Range("L210").Value2 = 6000
$ws.Range("M210").Value2 = 6000
$ws.Range("N210").Value2 = "`$/paquete 36 unidades"
$ws.Range("O210").Value2 = "Región Metropolitana"
$ws.Range("P210").Value2 = 167
$ws.Range("Q210").Value2 = 36
$ws.Range("R210").Value2 = "Hortaliza"
